# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 18:52"

# Helper to set multiple cell values for a row in one shot
function Set-RowValues {
    param(
        [int]$Row,
        [hashtable]$Values
    )
    foreach ($col in $Values.Keys) {
        $ws.Range("$col$Row").Value = $Values[$col]
    }
}

# Row 4 (Estados Unidos)
Set-RowValues 4 @{ B = 380744; C = 13740; E = 347528; G = 1036; H = 11907 }

# Row 7 (Alemania)
Set-RowValues 7 @{ B = 105604; C = 2229; E = 67618; G = 95; H = 1905 }

# Row 12 (Belgica)
Set-RowValues 12 @{ B = 34109; C = 3892; D = 1582; E = 31802; F = 1474; G = 76; H = 725 }

# Row 13 (Paises Bajos)
Set-RowValues 13 @{ E = 13365; G = 56; H = 821 }

# Row 19 (Israel)
Set-RowValues 19 @{ B = 12377; C = 194; E = 11668; G = 18; H = 582 }

# Row 21 (Noruega)
Set-RowValues 21 @{ E = 8262; G = 4; H = 61 }

# Row 27 (Pakistan)
Set-RowValues 27 @{ D = 421; E = 4601; G = 14; H = 150 }

# Row 33 (Arabia Saudita)
Set-RowValues 33 @{ B = 4009; C = 243; E = 3524; G = 3; H = 56 }

# Row 45 (Republica Dominicana)
Set-RowValues 45 @{ D = 888; E = 1343; F = 61 }

# Row 54 (Islandia)
Set-RowValues 54 @{ F = 76 }

# Row 83 (Kuwait)
Set-RowValues 83 @{ D = 62; E = 467; F = 29 }

# Row 110 (Mayotte)
Set-RowValues 110 @{ D = 46; E = 146 }

# Row 138 (San Pedro y Miquelon)
Set-RowValues 138 @{ E = 51; G = 1; H = 3 }
